# Update "想去人数" (F) and "最低票价" (G) figures that changed between
# crawler runs. Sheet 1 = "展览" (exhibitions), Sheet 4 = "全部类型"
# (all types) — the latter mirrors the same events (plus a few extra rows
# from the "演出"/"本地生活" sheets), so both need the same F/G updates.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)
$sheet4 = $wb.Worksheets.Item(4)

# --- Sheet 1 ("展览") ---
$sheet1.Range("F2").Value = 15234

$sheet1.Range("F4").Value = 82
$sheet1.Range("G4").Value = 55

$sheet1.Range("F5").Value = 1581
$sheet1.Range("G5").Value = 70

$sheet1.Range("F8").Value = 137

$sheet1.Range("F10").Value = 8413

$sheet1.Range("F12").Value = 55

$sheet1.Range("F15").Value = 1298

$sheet1.Range("F16").Value = 63

$sheet1.Range("F17").Value = 18

$sheet1.Range("F19").Value = 9026

$sheet1.Range("F21").Value = 86

$sheet1.Range("F22").Value = 202

$sheet1.Range("F23").Value = 162

$sheet1.Range("F24").Value = 328

$sheet1.Range("F25").Value = 5849

$sheet1.Range("F26").Value = 1030

$sheet1.Range("F27").Value = 40

$sheet1.Range("F29").Value = 82

# --- Sheet 4 ("全部类型") ---
$sheet4.Range("F2").Value = 15234

$sheet4.Range("F4").Value = 82
$sheet4.Range("G4").Value = 55

$sheet4.Range("F5").Value = 1581
$sheet4.Range("G5").Value = 70

$sheet4.Range("F9").Value = 137

$sheet4.Range("F11").Value = 8413

$sheet4.Range("F13").Value = 55

$sheet4.Range("F16").Value = 1298

$sheet4.Range("F17").Value = 63

$sheet4.Range("F18").Value = 18

$sheet4.Range("F22").Value = 9026

$sheet4.Range("F24").Value = 86

$sheet4.Range("F25").Value = 202

$sheet4.Range("F26").Value = 162

$sheet4.Range("F27").Value = 328

$sheet4.Range("F28").Value = 5849

$sheet4.Range("F29").Value = 1030

$sheet4.Range("F30").Value = 40

$sheet4.Range("F32").Value = 82
